# Automatische test-sync: 2025-06-19 19:04:30
# Adds the new incoming "Openingstijden" mail as row 35 on the "Logs"
# sheet, extends the conditional-formatting ranges to cover it, and
# re-sorts the category counts on the "Dashboard" sheet to reflect the
# bumped "Openingstijden" total.

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")
$dash = $wb.Worksheets.Item("Dashboard")

# --- 1. Append the new log entry on row 35 -------------------------------
$logs.Range("A35").Value = "Wat zijn jullie openingstijden?"
$logs.Range("B35").Value = "mailmind.test@zohomail.eu"
$logs.Range("C35").Value = "Hallo, ik zou graag willen weten wat jullie openingstijden zijn. Dank je wel!"
$logs.Range("D35").Value = "Openingstijden"
$logs.Range("E35").Value = "Beste klant,`nBedankt voor uw interesse. Onze openingstijden zijn van maandag tot en met vrijdag van 9:00 tot 18:00 uur. Op zaterdag zijn wij geopend van 10:00 tot 15:00 uur. Voor meer informatie of vragen kunt u altijd contact met ons opnemen.`nMet vriendelijke groet,`n[Naam bedrijf]"
$logs.Range("F35").Value = "2025-06-19 19:04:27"
$logs.Range("G35").Value = "Ja"

# --- 2. Extend the conditional-formatting sqrefs to include row 35 -------
# Use ModifyAppliesToRange on an existing rule so the dxf/style indices
# and rule order are preserved exactly (only the applied range changes).
$logs.Range("D2:D34").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D35"))
$logs.Range("G2:G34").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G35"))

# --- 3. Re-sort the Dashboard category table (rows 7-10) -----------------
# "Openingstijden" count goes 1 -> 2 (tying "Offerte-aanvraag"), so it
# now sorts ahead of it; the remaining ties keep their prior order.
$dash.Range("A7").Value = "Openingstijden"
$dash.Range("B7").Value = 2
$dash.Range("A8").Value = "Offerte-aanvraag"
$dash.Range("B8").Value = 2
$dash.Range("A9").Value = "Informatieaanvraag"
$dash.Range("B9").Value = 1
$dash.Range("A10").Value = "Samenwerking"
$dash.Range("B10").Value = 1
